$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2169893333333333
$ws.Range("H2").Value = 0.650968
$ws.Range("I2").Value = 0.2138091362408864
$ws.Range("J2").Value = 0.2138091362408864
$ws.Range("M2").Value = 281.0920463333333
$ws.Range("N2").Value = 843.2761389999999
$ws.Range("O2").Value = 0.8291026083535286
$ws.Range("P2").Value = 0.8291026083535286
$ws.Range("Q2").Value = 60.99397573917244
$ws.Range("R2").Value = 548.945781652552
$ws.Range("S2").Value = 0.1772697125471338
$ws.Range("T2").Value = 0.1772697125471338
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2169893333333333
$ws.Range("H3").Value = 0.650968
$ws.Range("I3").Value = 0.2138091362408864
$ws.Range("J3").Value = 0.2138091362408864
$ws.Range("O3").Value = 0.001324719879221983
$ws.Range("P3").Value = 0.001324719879221983
$ws.Range("Q3").Value = 0.09745468336533335
$ws.Range("R3").Value = 0.877092150288
$ws.Range("S3").Value = 0.0002832372131375836
$ws.Range("T3").Value = 0.0002832372131375836
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2169893333333333
$ws.Range("H4").Value = 0.650968
$ws.Range("I4").Value = 0.2138091362408864
$ws.Range("J4").Value = 0.2138091362408864
$ws.Range("M4").Value = 4.452417
$ws.Range("N4").Value = 13.357251
$ws.Range("O4").Value = 0.01313274635953239
$ws.Range("P4").Value = 0.01313274635953239
$ws.Range("Q4").Value = 0.966126996552
$ws.Range("R4").Value = 8.695142968968
$ws.Range("S4").Value = 0.002807901155602265
$ws.Range("T4").Value = 0.002807901155602265
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2169893333333333
$ws.Range("H5").Value = 0.650968
$ws.Range("I5").Value = 0.2138091362408864
$ws.Range("J5").Value = 0.2138091362408864
$ws.Range("M5").Value = 53.03808999999999
$ws.Range("N5").Value = 159.11427
$ws.Range("O5").Value = 0.156439925407717
$ws.Range("P5").Value = 0.156439925407717
$ws.Range("Q5").Value = 11.50869979037333
$ws.Range("R5").Value = 103.57829811336
$ws.Range("S5").Value = 0.03344828532501266
$ws.Range("T5").Value = 0.03344828532501267
$ws.Range("I6").Value = 0.3272432505578689
$ws.Range("J6").Value = 0.3272432505578689
$ws.Range("M6").Value = 281.0920463333333
$ws.Range("N6").Value = 843.2761389999999
$ws.Range("O6").Value = 0.8291026083535286
$ws.Range("P6").Value = 0.8291026083535286
$ws.Range("Q6").Value = 93.35366690246089
$ws.Range("R6").Value = 840.183002122148
$ws.Range("S6").Value = 0.2713182326036164
$ws.Range("T6").Value = 0.2713182326036164
$ws.Range("I7").Value = 0.3272432505578689
$ws.Range("J7").Value = 0.3272432505578689
$ws.Range("O7").Value = 0.001324719879221983
$ws.Range("P7").Value = 0.001324719879221983
$ws.Range("S7").Value = 0.0004335056393552293
$ws.Range("T7").Value = 0.0004335056393552293
$ws.Range("I8").Value = 0.3272432505578689
$ws.Range("J8").Value = 0.3272432505578689
$ws.Range("M8").Value = 4.452417
$ws.Range("N8").Value = 13.357251
$ws.Range("O8").Value = 0.01313274635953239
$ws.Range("P8").Value = 0.01313274635953239
$ws.Range("Q8").Value = 1.478695178148
$ws.Range("R8").Value = 13.308256603332
$ws.Range("S8").Value = 0.004297602607445398
$ws.Range("T8").Value = 0.004297602607445398
$ws.Range("I9").Value = 0.3272432505578689
$ws.Range("J9").Value = 0.3272432505578689
$ws.Range("M9").Value = 53.03808999999999
$ws.Range("N9").Value = 159.11427
$ws.Range("O9").Value = 0.156439925407717
$ws.Range("P9").Value = 0.156439925407717
$ws.Range("Q9").Value = 17.61451542862666
$ws.Range("R9").Value = 158.53063885764
$ws.Range("S9").Value = 0.05119390970745185
$ws.Range("T9").Value = 0.05119390970745186
$ws.Range("G10").Value = 0.4136746666666666
$ws.Range("H10").Value = 1.241024
$ws.Range("I10").Value = 0.4076118480389355
$ws.Range("J10").Value = 0.4076118480389355
$ws.Range("M10").Value = 281.0920463333333
$ws.Range("N10").Value = 843.2761389999999
$ws.Range("O10").Value = 0.8291026083535286
$ws.Range("P10").Value = 0.8291026083535286
$ws.Range("Q10").Value = 116.2806585695929
$ws.Range("R10").Value = 1046.525927126336
$ws.Range("S10").Value = 0.3379520464048835
$ws.Range("T10").Value = 0.3379520464048836
$ws.Range("G11").Value = 0.4136746666666666
$ws.Range("H11").Value = 1.241024
$ws.Range("I11").Value = 0.4076118480389355
$ws.Range("J11").Value = 0.4076118480389355
$ws.Range("O11").Value = 0.001324719879221983
$ws.Range("P11").Value = 0.001324719879221983
$ws.Range("Q11").Value = 0.1857903936426667
$ws.Range("R11").Value = 1.672113542784
$ws.Range("S11").Value = 0.000539971518103588
$ws.Range("T11").Value = 0.0005399715181035881
$ws.Range("G12").Value = 0.4136746666666666
$ws.Range("H12").Value = 1.241024
$ws.Range("I12").Value = 0.4076118480389355
$ws.Range("J12").Value = 0.4076118480389355
$ws.Range("M12").Value = 4.452417
$ws.Range("N12").Value = 13.357251
$ws.Range("O12").Value = 0.01313274635953239
$ws.Range("P12").Value = 0.01313274635953239
$ws.Range("Q12").Value = 1.841852118336
$ws.Range("R12").Value = 16.576669065024
$ws.Range("S12").Value = 0.005353063013435598
$ws.Range("T12").Value = 0.0053530630134356
$ws.Range("G13").Value = 0.4136746666666666
$ws.Range("H13").Value = 1.241024
$ws.Range("I13").Value = 0.4076118480389355
$ws.Range("J13").Value = 0.4076118480389355
$ws.Range("M13").Value = 53.03808999999999
$ws.Range("N13").Value = 159.11427
$ws.Range("O13").Value = 0.156439925407717
$ws.Range("P13").Value = 0.156439925407717
$ws.Range("Q13").Value = 21.94051420138666
$ws.Range("R13").Value = 197.46462781248
$ws.Range("S13").Value = 0.06376676710251274
$ws.Range("T13").Value = 0.06376676710251275
$ws.Range("G14").Value = 0.05209933333333334
$ws.Range("H14").Value = 0.156298
$ws.Range("I14").Value = 0.05133576516230915
$ws.Range("J14").Value = 0.05133576516230916
$ws.Range("M14").Value = 281.0920463333333
$ws.Range("N14").Value = 843.2761389999999
$ws.Range("O14").Value = 0.8291026083535286
$ws.Range("P14").Value = 0.8291026083535286
$ws.Range("Q14").Value = 14.64470821926911
$ws.Range("R14").Value = 131.802373973422
$ws.Range("S14").Value = 0.04256261679789472
$ws.Range("T14").Value = 0.04256261679789473
$ws.Range("G15").Value = 0.05209933333333334
$ws.Range("H15").Value = 0.156298
$ws.Range("I15").Value = 0.05133576516230915
$ws.Range("J15").Value = 0.05133576516230916
$ws.Range("O15").Value = 0.001324719879221983
$ws.Range("P15").Value = 0.001324719879221983
$ws.Range("Q15").Value = 0.02339895678533334
$ws.Range("R15").Value = 0.210590611068
$ws.Range("S15").Value = 0.00006800550862558228
$ws.Range("T15").Value = 0.0000680055086255823
$ws.Range("G16").Value = 0.05209933333333334
$ws.Range("H16").Value = 0.156298
$ws.Range("I16").Value = 0.05133576516230915
$ws.Range("J16").Value = 0.05133576516230916
$ws.Range("M16").Value = 4.452417
$ws.Range("N16").Value = 13.357251
$ws.Range("O16").Value = 0.01313274635953239
$ws.Range("P16").Value = 0.01313274635953239
$ws.Range("Q16").Value = 0.231967957422
$ws.Range("R16").Value = 2.087711616798
$ws.Range("S16").Value = 0.0006741795830491251
$ws.Range("T16").Value = 0.0006741795830491253
$ws.Range("G17").Value = 0.05209933333333334
$ws.Range("H17").Value = 0.156298
$ws.Range("I17").Value = 0.05133576516230915
$ws.Range("J17").Value = 0.05133576516230916
$ws.Range("M17").Value = 53.03808999999999
$ws.Range("N17").Value = 159.11427
$ws.Range("O17").Value = 0.156439925407717
$ws.Range("P17").Value = 0.156439925407717
$ws.Range("Q17").Value = 2.763249130273333
$ws.Range("R17").Value = 24.86924217246
$ws.Range("S17").Value = 0.00803096327273972
$ws.Range("T17").Value = 0.008030963272739722
